# Updated cryptos list - applies price/volume changes from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Numeric-looking strings (single decimal point) must be forced to
    # Text format first, otherwise Excel auto-converts them to numbers
    # and silently drops the original text formatting (e.g. trailing zeros).
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$ws.Range("D2").Value = "28.488.65"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.957.73"
$ws.Range("E3").Value = "  -0.29%  "
Set-TextValue $ws.Range("D4") "1.012"
$ws.Range("E4").Value = "  +0.61%  "
Set-TextValue $ws.Range("D5") "322.63"
$ws.Range("E5").Value = "  -1.34%  "
Set-TextValue $ws.Range("D6") "1.010"
$ws.Range("E6").Value = "  +0.57%  "
Set-TextValue $ws.Range("D7") "0.4798"
$ws.Range("E7").Value = "  -3.86%  "
Set-TextValue $ws.Range("D8") "0.4063"
$ws.Range("E8").Value = "  -3.50%  "
Set-TextValue $ws.Range("D9") "54.07"
$ws.Range("E9").Value = "  +1.52%  "
Set-TextValue $ws.Range("D10") "0.08504"
$ws.Range("E10").Value = "  -7.87%  "
Set-TextValue $ws.Range("D11") "1.059"
$ws.Range("E11").Value = "  -3.69%  "
Set-TextValue $ws.Range("D12") "22.39"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "1.976.68"
$ws.Range("E13").Value = "  +0.63%  "
Set-TextValue $ws.Range("D14") "7.575"
$ws.Range("E14").Value = "  -3.72%  "
Set-TextValue $ws.Range("D15") "6.171"
$ws.Range("E15").Value = "  -4.32%  "
Set-TextValue $ws.Range("D16") "1.012"
$ws.Range("E16").Value = "  +0.67%  "
Set-TextValue $ws.Range("D17") "90.74"
$ws.Range("E17").Value = "  -0.71%  "
Set-TextValue $ws.Range("D18") "0.00001072"
$ws.Range("E18").Value = "  -2.61%  "
Set-TextValue $ws.Range("D19") "0.06621"
$ws.Range("E19").Value = "  -1.11%  "
Set-TextValue $ws.Range("D20") "18.45"
$ws.Range("E20").Value = "  -4.01%  "
Set-TextValue $ws.Range("D21") "1.010"
$ws.Range("E21").Value = "  +0.57%  "
Set-TextValue $ws.Range("D22") "5.849"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "28.509.25"
$ws.Range("E23").Value = "  -1.84%  "
Set-TextValue $ws.Range("D24") "11.44"
$ws.Range("E24").Value = "  -5.16%  "
Set-TextValue $ws.Range("D25") "2.295"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "2.209.12"
$ws.Range("E26").Value = "  +0.28%  "
Set-TextValue $ws.Range("D27") "156.00"
$ws.Range("E27").Value = "  -0.36%  "
Set-TextValue $ws.Range("D28") "20.30"
$ws.Range("E28").Value = "  -1.44%  "
Set-TextValue $ws.Range("D29") "2.179"
$ws.Range("E29").Value = "  -3.72%  "
Set-TextValue $ws.Range("D30") "5.856"
$ws.Range("E30").Value = "  -5.80%  "
Set-TextValue $ws.Range("D31") "124.48"
$ws.Range("E31").Value = "  -1.95%  "
Set-TextValue $ws.Range("D32") "0.9872"
$ws.Range("E32").Value = "  -5.33%  "
Set-TextValue $ws.Range("D33") "0.09667"
$ws.Range("E33").Value = "  -1.89%  "
Set-TextValue $ws.Range("D34") "1.456"
$ws.Range("E34").Value = "  -4.76%  "
Set-TextValue $ws.Range("D37") "9.110"
$ws.Range("E37").Value = "  +2.03%  "
Set-TextValue $ws.Range("D38") "0.02336"
$ws.Range("E38").Value = "  -3.59%  "
Set-TextValue $ws.Range("D39") "0.06214"
$ws.Range("E39").Value = "  -1.63%  "
Set-TextValue $ws.Range("D40") "1.253"
$ws.Range("E40").Value = "  -3.67%  "
Set-TextValue $ws.Range("D41") "0.6225"
$ws.Range("E41").Value = "  -3.51%  "
Set-TextValue $ws.Range("D42") "11.20"
$ws.Range("E42").Value = "  -2.08%  "
Set-TextValue $ws.Range("D43") "1.010"
$ws.Range("E43").Value = "  +0.76%  "
Set-TextValue $ws.Range("D44") "0.1916"
$ws.Range("E44").Value = "  -3.69%  "
Set-TextValue $ws.Range("D45") "1.358"
$ws.Range("E45").Value = "  +5.90%  "
Set-TextValue $ws.Range("D46") "0.5955"
$ws.Range("E46").Value = "  -4.46%  "
Set-TextValue $ws.Range("D47") "13.01"
$ws.Range("E47").Value = "  -3.02%  "
Set-TextValue $ws.Range("D48") "2.060"
$ws.Range("E48").Value = "  -5.49%  "
Set-TextValue $ws.Range("D49") "3.411"
$ws.Range("E49").Value = "  -1.40%  "
Set-TextValue $ws.Range("D50") "0.06832"
$ws.Range("E50").Value = "  -1.23%  "
Set-TextValue $ws.Range("D51") "111.29"
$ws.Range("E51").Value = "  -1.13%  "

# Rows 35 and 36 swap their entire content (HuobiToken <-> Filecoin)
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D35") "5.638"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D36") "3.691"
$ws.Range("E36").Value = "  +0.52%  "
